# Apply updated cryptocurrency price / volume(1h) figures to sheet1.
# Values are written with a leading apostrophe to force text storage
# (matching the source data, which is all inline/shared text, not numbers),
# then the cell style is reset to "Normal" so no stray quote-prefix /
# number-format style gets introduced by Excel's auto-detection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.295.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.65%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.665.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.57%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.85%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'219.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.58%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +1.39%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.80%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2647"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.37%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.32%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.60%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07822"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.32%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.566"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.12%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.667.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.68%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.892.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.64%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0₅8208"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.48%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'65.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.24%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +0.86%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.707"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.29%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'193.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.00%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'6.042"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.09%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.84%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'145.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.66%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.1233"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.49%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.202"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.74%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'16.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.03%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +3.73%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.05900"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.16%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.281"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.29%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.625"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.07%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.279"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.38%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.610"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.85%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.9645"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.22%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.825"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.51%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.418"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.29%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.5808"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.72%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01610"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.67%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.8657"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.95%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.868"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.79%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.051.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.00%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.010"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.77%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'104.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.57%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.803.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.30%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'57.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "'  -4.97%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.013"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.14%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +1.86%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'8.029"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.33%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.05167"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.28%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.417"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.18%  "
$ws.Range("E51").Style = "Normal"
